# Apply updated PSSM values (supplemental figures) to Sheet1, cells B2:K21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = -18.83914841613093
$ws.Cells.Item(2, 3).Value = -18.83914841613093
$ws.Cells.Item(2, 4).Value = -18.83914841613093
$ws.Cells.Item(2, 5).Value = -18.83914841613093
$ws.Cells.Item(2, 6).Value = -18.83914841613093
$ws.Cells.Item(2, 7).Value = -18.83914841613093
$ws.Cells.Item(2, 8).Value = -18.83914841613093
$ws.Cells.Item(2, 9).Value = -18.83914841613093
$ws.Cells.Item(2, 10).Value = -18.83914841613093
$ws.Cells.Item(2, 11).Value = -18.83914841613093

$ws.Cells.Item(3, 2).Value = -18.83914841613093
$ws.Cells.Item(3, 3).Value = -18.83914841613093
$ws.Cells.Item(3, 4).Value = -18.83914841613093
$ws.Cells.Item(3, 5).Value = -18.83914841613093
$ws.Cells.Item(3, 6).Value = -18.83914841613093
$ws.Cells.Item(3, 7).Value = -18.83914841613093
$ws.Cells.Item(3, 8).Value = -18.83914841613093
$ws.Cells.Item(3, 9).Value = 2.402679410809143
$ws.Cells.Item(3, 10).Value = -18.83914841613093
$ws.Cells.Item(3, 11).Value = -18.83914841613093

$ws.Cells.Item(4, 2).Value = -18.83914841613093
$ws.Cells.Item(4, 3).Value = -18.83914841613093
$ws.Cells.Item(4, 4).Value = 3.058319294891543
$ws.Cells.Item(4, 5).Value = -18.83914841613093
$ws.Cells.Item(4, 6).Value = 2.502479530253949
$ws.Cells.Item(4, 7).Value = -18.83914841613093
$ws.Cells.Item(4, 8).Value = 1.852735643661177
$ws.Cells.Item(4, 9).Value = -18.83914841613093
$ws.Cells.Item(4, 10).Value = 2.346059137730396
$ws.Cells.Item(4, 11).Value = -18.83914841613093

$ws.Cells.Item(5, 2).Value = -18.83914841613093
$ws.Cells.Item(5, 3).Value = -18.83914841613093
$ws.Cells.Item(5, 4).Value = -18.83914841613093
$ws.Cells.Item(5, 5).Value = -18.83914841613093
$ws.Cells.Item(5, 6).Value = -18.83914841613093
$ws.Cells.Item(5, 7).Value = 2.170347415295987
$ws.Cells.Item(5, 8).Value = -18.83914841613093
$ws.Cells.Item(5, 9).Value = -18.83914841613093
$ws.Cells.Item(5, 10).Value = -18.83914841613093
$ws.Cells.Item(5, 11).Value = -18.83914841613093

$ws.Cells.Item(6, 2).Value = -18.83914841613093
$ws.Cells.Item(6, 3).Value = -18.83914841613093
$ws.Cells.Item(6, 4).Value = -18.83914841613093
$ws.Cells.Item(6, 5).Value = -18.83914841613093
$ws.Cells.Item(6, 6).Value = -18.83914841613093
$ws.Cells.Item(6, 7).Value = -18.83914841613093
$ws.Cells.Item(6, 8).Value = -18.83914841613093
$ws.Cells.Item(6, 9).Value = -18.83914841613093
$ws.Cells.Item(6, 10).Value = -18.83914841613093
$ws.Cells.Item(6, 11).Value = -18.83914841613093

$ws.Cells.Item(7, 2).Value = 2.998218118827072
$ws.Cells.Item(7, 3).Value = -18.83914841613093
$ws.Cells.Item(7, 4).Value = -18.83914841613093
$ws.Cells.Item(7, 5).Value = -18.83914841613093
$ws.Cells.Item(7, 6).Value = -18.83914841613093
$ws.Cells.Item(7, 7).Value = -18.83914841613093
$ws.Cells.Item(7, 8).Value = -18.83914841613093
$ws.Cells.Item(7, 9).Value = -18.83914841613093
$ws.Cells.Item(7, 10).Value = -18.83914841613093
$ws.Cells.Item(7, 11).Value = -18.83914841613093

$ws.Cells.Item(8, 2).Value = -18.83914841613093
$ws.Cells.Item(8, 3).Value = -18.83914841613093
$ws.Cells.Item(8, 4).Value = -18.83914841613093
$ws.Cells.Item(8, 5).Value = 3.004119787294508
$ws.Cells.Item(8, 6).Value = -18.83914841613093
$ws.Cells.Item(8, 7).Value = -18.83914841613093
$ws.Cells.Item(8, 8).Value = -18.83914841613093
$ws.Cells.Item(8, 9).Value = -18.83914841613093
$ws.Cells.Item(8, 10).Value = -18.83914841613093
$ws.Cells.Item(8, 11).Value = -18.83914841613093

$ws.Cells.Item(9, 2).Value = 3.58614458937079
$ws.Cells.Item(9, 3).Value = -18.83914841613093
$ws.Cells.Item(9, 4).Value = -18.83914841613093
$ws.Cells.Item(9, 5).Value = -18.83914841613093
$ws.Cells.Item(9, 6).Value = -18.83914841613093
$ws.Cells.Item(9, 7).Value = -18.83914841613093
$ws.Cells.Item(9, 8).Value = -18.83914841613093
$ws.Cells.Item(9, 9).Value = -18.83914841613093
$ws.Cells.Item(9, 10).Value = -18.83914841613093
$ws.Cells.Item(9, 11).Value = -18.83914841613093

$ws.Cells.Item(10, 2).Value = -18.83914841613093
$ws.Cells.Item(10, 3).Value = -18.83914841613093
$ws.Cells.Item(10, 4).Value = -18.83914841613093
$ws.Cells.Item(10, 5).Value = -18.83914841613093
$ws.Cells.Item(10, 6).Value = -18.83914841613093
$ws.Cells.Item(10, 7).Value = -18.83914841613093
$ws.Cells.Item(10, 8).Value = -18.83914841613093
$ws.Cells.Item(10, 9).Value = 1.58463330127254
$ws.Cells.Item(10, 10).Value = -18.83914841613093
$ws.Cells.Item(10, 11).Value = 2.212600207507286

$ws.Cells.Item(11, 2).Value = -18.83914841613093
$ws.Cells.Item(11, 3).Value = -18.83914841613093
$ws.Cells.Item(11, 4).Value = -18.83914841613093
$ws.Cells.Item(11, 5).Value = 2.201714510935344
$ws.Cells.Item(11, 6).Value = -18.83914841613093
$ws.Cells.Item(11, 7).Value = 2.626312866717237
$ws.Cells.Item(11, 8).Value = -18.83914841613093
$ws.Cells.Item(11, 9).Value = -18.83914841613093
$ws.Cells.Item(11, 10).Value = -18.83914841613093
$ws.Cells.Item(11, 11).Value = 1.371428999138703

$ws.Cells.Item(12, 2).Value = -18.83914841613093
$ws.Cells.Item(12, 3).Value = -18.83914841613093
$ws.Cells.Item(12, 4).Value = -18.83914841613093
$ws.Cells.Item(12, 5).Value = -18.83914841613093
$ws.Cells.Item(12, 6).Value = -18.83914841613093
$ws.Cells.Item(12, 7).Value = -18.83914841613093
$ws.Cells.Item(12, 8).Value = -18.83914841613093
$ws.Cells.Item(12, 9).Value = -18.83914841613093
$ws.Cells.Item(12, 10).Value = -18.83914841613093
$ws.Cells.Item(12, 11).Value = -18.83914841613093

$ws.Cells.Item(13, 2).Value = -18.83914841613093
$ws.Cells.Item(13, 3).Value = -18.83914841613093
$ws.Cells.Item(13, 4).Value = -18.83914841613093
$ws.Cells.Item(13, 5).Value = 1.851048251942255
$ws.Cells.Item(13, 6).Value = -18.83914841613093
$ws.Cells.Item(13, 7).Value = -18.83914841613093
$ws.Cells.Item(13, 8).Value = -18.83914841613093
$ws.Cells.Item(13, 9).Value = -18.83914841613093
$ws.Cells.Item(13, 10).Value = 2.273001521730748
$ws.Cells.Item(13, 11).Value = 1.607568413899657

$ws.Cells.Item(14, 2).Value = -18.83914841613093
$ws.Cells.Item(14, 3).Value = -18.83914841613093
$ws.Cells.Item(14, 4).Value = 1.959653220223898
$ws.Cells.Item(14, 5).Value = -18.83914841613093
$ws.Cells.Item(14, 6).Value = -18.83914841613093
$ws.Cells.Item(14, 7).Value = -18.83914841613093
$ws.Cells.Item(14, 8).Value = -18.83914841613093
$ws.Cells.Item(14, 9).Value = -18.83914841613093
$ws.Cells.Item(14, 10).Value = -18.83914841613093
$ws.Cells.Item(14, 11).Value = 2.139666064695753

$ws.Cells.Item(15, 2).Value = -18.83914841613093
$ws.Cells.Item(15, 3).Value = -18.83914841613093
$ws.Cells.Item(15, 4).Value = -0.358527803506565
$ws.Cells.Item(15, 5).Value = -18.83914841613093
$ws.Cells.Item(15, 6).Value = -18.83914841613093
$ws.Cells.Item(15, 7).Value = -18.83914841613093
$ws.Cells.Item(15, 8).Value = -18.83914841613093
$ws.Cells.Item(15, 9).Value = -18.83914841613093
$ws.Cells.Item(15, 10).Value = -18.83914841613093
$ws.Cells.Item(15, 11).Value = -18.83914841613093

$ws.Cells.Item(16, 2).Value = -18.83914841613093
$ws.Cells.Item(16, 3).Value = -18.83914841613093
$ws.Cells.Item(16, 4).Value = -18.83914841613093
$ws.Cells.Item(16, 5).Value = -18.83914841613093
$ws.Cells.Item(16, 6).Value = -18.83914841613093
$ws.Cells.Item(16, 7).Value = -18.83914841613093
$ws.Cells.Item(16, 8).Value = -18.83914841613093
$ws.Cells.Item(16, 9).Value = -18.83914841613093
$ws.Cells.Item(16, 10).Value = 2.302884454776093
$ws.Cells.Item(16, 11).Value = -18.83914841613093

$ws.Cells.Item(17, 2).Value = -18.83914841613093
$ws.Cells.Item(17, 3).Value = -18.83914841613093
$ws.Cells.Item(17, 4).Value = 0.9927417488227375
$ws.Cells.Item(17, 5).Value = -18.83914841613093
$ws.Cells.Item(17, 6).Value = -18.83914841613093
$ws.Cells.Item(17, 7).Value = -18.83914841613093
$ws.Cells.Item(17, 8).Value = 0.6685475215490598
$ws.Cells.Item(17, 9).Value = 1.011235315729628
$ws.Cells.Item(17, 10).Value = 1.271142441986408
$ws.Cells.Item(17, 11).Value = -18.83914841613093

$ws.Cells.Item(18, 2).Value = -18.83914841613093
$ws.Cells.Item(18, 3).Value = -18.83914841613093
$ws.Cells.Item(18, 4).Value = -18.83914841613093
$ws.Cells.Item(18, 5).Value = -18.83914841613093
$ws.Cells.Item(18, 6).Value = -18.83914841613093
$ws.Cells.Item(18, 7).Value = -18.83914841613093
$ws.Cells.Item(18, 8).Value = 0.6757678556830267
$ws.Cells.Item(18, 9).Value = 1.185711235551648
$ws.Cells.Item(18, 10).Value = 1.451222608281994
$ws.Cells.Item(18, 11).Value = -18.83914841613093

$ws.Cells.Item(19, 2).Value = -18.83914841613093
$ws.Cells.Item(19, 3).Value = -18.83914841613093
$ws.Cells.Item(19, 4).Value = 1.057877174421247
$ws.Cells.Item(19, 5).Value = -18.83914841613093
$ws.Cells.Item(19, 6).Value = -18.83914841613093
$ws.Cells.Item(19, 7).Value = -18.83914841613093
$ws.Cells.Item(19, 8).Value = 1.84720591505243
$ws.Cells.Item(19, 9).Value = 2.000065369455606
$ws.Cells.Item(19, 10).Value = -18.83914841613093
$ws.Cells.Item(19, 11).Value = -18.83914841613093

$ws.Cells.Item(20, 2).Value = -18.83914841613093
$ws.Cells.Item(20, 3).Value = -18.83914841613093
$ws.Cells.Item(20, 4).Value = 1.550152063135593
$ws.Cells.Item(20, 5).Value = -18.83914841613093
$ws.Cells.Item(20, 6).Value = 3.841306617998562
$ws.Cells.Item(20, 7).Value = -18.83914841613093
$ws.Cells.Item(20, 8).Value = 2.180250240881787
$ws.Cells.Item(20, 9).Value = 1.774940415916563
$ws.Cells.Item(20, 10).Value = -18.83914841613093
$ws.Cells.Item(20, 11).Value = 2.412408539635536

$ws.Cells.Item(21, 2).Value = -18.83914841613093
$ws.Cells.Item(21, 3).Value = 4.321925172419001
$ws.Cells.Item(21, 4).Value = -18.83914841613093
$ws.Cells.Item(21, 5).Value = 1.914261776048808
$ws.Cells.Item(21, 6).Value = -18.83914841613093
$ws.Cells.Item(21, 7).Value = 3.220977848091739
$ws.Cells.Item(21, 8).Value = 2.34223014208036
$ws.Cells.Item(21, 9).Value = -18.83914841613093
$ws.Cells.Item(21, 10).Value = -18.83914841613093
$ws.Cells.Item(21, 11).Value = -18.83914841613093
